# Add 2022-Q4 data
# 1) Insert a new worksheet named "2022-Q4" right after "总计" and before "2022-Q2".
# 2) Populate it with the fund holding data for 2022-Q4.
# 3) Update the "总计" (summary) sheet to add a new row for 2022-Q4 on top of
#    the existing quarters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create + position the new "2022-Q4" sheet
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add()
$q4.Name = "2022-Q4"
$q2 = $wb.Worksheets.Item("2022-Q2")
$q4.Move($q2)
# Re-fetch the worksheet reference: the old $q4 COM reference can go stale
# once the sheet has been repositioned.
$q4 = $wb.Worksheets.Item("2022-Q4")

# ---------------------------------------------------------------------------
# Helper to style a header-ish cell the same way the rest of the workbook
# styles its "label" cells: bold, centered, top aligned, thin box border.
# ---------------------------------------------------------------------------
function Set-LabelStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# ---------------------------------------------------------------------------
# Step 2: fill in the "2022-Q4" sheet contents
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    Set-LabelStyle $cell
}

$q4Rows = @(
    @(0, "001037", "国投瑞银锐意改革灵活配置混合A", "1.99", "85.61", "4.23", "0.0842", 6),
    @(1, "016780", "国投瑞银锐意改革灵活配置混合C", "0.00", "85.61", "4.23", 0,        6)
)

$r = 2
foreach ($row in $q4Rows) {
    $aCell = $q4.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Set-LabelStyle $aCell

    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]

    $gCell = $q4.Cells.Item($r, 7)
    if ($row[6] -is [string]) {
        $gCell.Value = "'" + $row[6]
    } else {
        $gCell.Value = $row[6]
    }

    $q4.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# ---------------------------------------------------------------------------
# Step 3: update the "总计" sheet with the new 2022-Q4 row
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Cells.Clear()

$sumHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $sumHeaders.Length; $i++) {
    $cell = $summary.Cells.Item(1, $i + 2)
    $cell.Value = $sumHeaders[$i]
    Set-LabelStyle $cell
}

$sumRows = @(
    @(0, "2022-Q4", 2,  0.08),
    @(1, "2022-Q2", 7,  0.92),
    @(2, "2022-Q1", 4,  0.6),
    @(3, "2021-Q4", 10, 1.54)
)

$r = 2
foreach ($row in $sumRows) {
    $aCell = $summary.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Set-LabelStyle $aCell

    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]

    $r++
}
